# Commit: swap the presentation's two themes -- ppt/theme/theme1.xml
# (the "Integral" theme used by the slide master / all slides) and
# ppt/theme/theme2.xml (the "Office Theme" used by the notes master)
# traded their DrawingML colour schemes (fontScheme/fmtScheme were
# already identical between the two parts).
#
# The PowerPoint object model exposes the 12-colour DrawingML theme
# scheme via Slide.ThemeColorScheme (any slide resolves to the single
# slide master's theme -- ppt/theme/theme1.xml in this deck). The
# ThemeColorScheme.Colors(index) collection is ordered exactly like
# MsoThemeColorSchemeIndex:
#   1 = Dark1, 2 = Light1, 3 = Dark2, 4 = Light2,
#   5..10 = Accent1..Accent6, 11 = Hyperlink, 12 = FollowedHyperlink
# and RGBColor.RGB uses the usual COM 0xBBGGRR byte order.

$p = $ppt.ActivePresentation

# Target palette: the "Office Theme" colour scheme that used to live in
# ppt/theme/theme2.xml, now applied to the main theme (theme1.xml).
$officeThemeColors = @(
    0x000000,  #  1 dk1      srgbClr 000000
    0xFFFFFF,  #  2 lt1      srgbClr FFFFFF
    0x6A5444,  #  3 dk2      srgbClr 44546A
    0xE6E6E7,  #  4 lt2      srgbClr E7E6E6
    0xD59B5B,  #  5 accent1  srgbClr 5B9BD5
    0x317DED,  #  6 accent2  srgbClr ED7D31
    0xA5A5A5,  #  7 accent3  srgbClr A5A5A5
    0x00C0FF,  #  8 accent4  srgbClr FFC000
    0xC47244,  #  9 accent5  srgbClr 4472C4
    0x47AD70,  # 10 accent6  srgbClr 70AD47
    0xC16305,  # 11 hlink    srgbClr 0563C1
    0x724F95   # 12 folHlink srgbClr 954F72
)

$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
